$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, shifting existing rows 49-63 down to 50-64.
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with the new weekly price record.
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44704
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100108
$ws.Range("H49").Value = "Tropicales y subtropicales"
$ws.Range("I49").Value = 100108007
$ws.Range("J49").Value = "Coco"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 12
$ws.Range("N49").Value = 30000
$ws.Range("O49").Value = 30000
$ws.Range("P49").Value = 30000
$ws.Range("Q49").Value = "$/malla 20 unidades"
$ws.Range("R49").Value = "Perú"
$ws.Range("S49").Value = 1500
$ws.Range("T49").Value = 20
